# "Added anaysis UK Pillar 2 SGTF data"
#
# Nudges the position/size of a handful of scatter-point marker shapes
# (and resizes/relabels two of their numeric value-labels) on slide 1.
#
# All affected shapes live inside the single top-level group shape
# (Shapes.Item(1)); within that group the GroupItems index is simply
# (shape-id - 2), since the group's direct children are numbered
# contiguously starting at id 3.
#
# Shape.Left/Top/Width/Height round-trip through a 32-bit float internally,
# so naively writing `emu / 12700` can land one EMU short of the exact
# target after it gets truncated back down. EmuToPointsPrecise nudges the
# point value by progressively larger steps until the float32 round-trip
# reproduces the exact target EMU, so the saved XML matches precisely.

$EMU_PER_POINT = 12700

function EmuToPointsPrecise {
    param($TargetEmu)
    $pts = [double]$TargetEmu / $EMU_PER_POINT
    $step = 0.0000001
    for ($i = 0; $i -lt 60; $i++) {
        $f = [float]$pts
        $achieved = [math]::Floor([double]$f * $EMU_PER_POINT)
        if ($achieved -eq $TargetEmu) {
            break
        }
        if ($achieved -lt $TargetEmu) {
            $pts = $pts + $step
        } else {
            $pts = $pts - $step
        }
        $step = $step * 1.5
    }
    return $pts
}

function Set-ShapeGeometryEmu {
    param($Shape, $LeftEmu, $TopEmu, $WidthEmu, $HeightEmu)
    $Shape.Left   = EmuToPointsPrecise $LeftEmu
    $Shape.Top    = EmuToPointsPrecise $TopEmu
    $Shape.Width  = EmuToPointsPrecise $WidthEmu
    $Shape.Height = EmuToPointsPrecise $HeightEmu
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

# pt14
Set-ShapeGeometryEmu $g.GroupItems.Item(12) 2181658 6070370 67453 67453

# pt15
Set-ShapeGeometryEmu $g.GroupItems.Item(13) 2771663 6075092 58010 58010

# pt16
Set-ShapeGeometryEmu $g.GroupItems.Item(14) 3352225 4923277 67453 67453

# pt17
Set-ShapeGeometryEmu $g.GroupItems.Item(15) 3937508 3488708 67453 67453

# pt18
Set-ShapeGeometryEmu $g.GroupItems.Item(16) 4509990 3426456 93055 93055

# pt19
Set-ShapeGeometryEmu $g.GroupItems.Item(17) 5091277 3399324 101049 101049

# pt20
Set-ShapeGeometryEmu $g.GroupItems.Item(18) 5671298 2477350 111573 111573

# pt47
Set-ShapeGeometryEmu $g.GroupItems.Item(45) 7629570 3578317 9271 9271

# pt49
Set-ShapeGeometryEmu $g.GroupItems.Item(47) 7603069 3771272 62273 62273

# tx52 (label box resized + text "10" -> "1")
$tx52 = $g.GroupItems.Item(50)
Set-ShapeGeometryEmu $tx52 7819849 3538959 67806 87630
$tx52.TextFrame.TextRange.Text = "1"

# tx53 (label text "50" -> "10")
$tx53 = $g.GroupItems.Item(51)
$tx53.TextFrame.TextRange.Text = "10"
